# Scheduled runner update: refresh market-derived price/profit columns
# (H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#  K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ)
# across several sheets/rows.

$wb = $excel.ActiveWorkbook

function Set-Row {
    param(
        $ws,
        [int]$row,
        [double[]]$values   # values for H,I,J,K,L,M,N in order; $null entries are skipped
    )
    $cols = 8,9,10,11,12,13,14
    for ($i = 0; $i -lt $cols.Length; $i++) {
        if ($null -ne $values[$i]) {
            $ws.Cells.Item($row, $cols[$i]).Value = $values[$i]
        }
    }
}

function Clear-Cell {
    param($ws, [int]$row, [int]$col)
    $ws.Cells.Item($row, $col).ClearContents()
}

# ---------------- ALC ----------------
$ws = $wb.Worksheets.Item("ALC")

Set-Row $ws 32  @(628.8, 817.5, 413.14285, 817.5, 413.14285, -491.5, -1065.14285)
Set-Row $ws 125 @(1351.75, 718.4, 1639.6364, 6465.599999999999, 14756.7276, -4005.599999999999, -19676.7276)
Set-Row $ws 131 @(1315.5, 404.8, 2833.3333, 1214.4, 8499.999899999999, 3825.6, -18579.9999)
Set-Row $ws 132 @(3140.35, 2927.9421, 4472.727, 8783.826300000001, 13418.181, -6253.826300000001, -18478.181)

# ---------------- ARM ----------------
$ws = $wb.Worksheets.Item("ARM")

Set-Row $ws 32  @(47016.293, 44956.6, 48100.344, 44956.6, 48100.344, -44669.6, -48674.344)

Set-Row $ws 44  @(0, 0, 0, 0, 0, $null, $null)
Clear-Cell $ws 44 14

Set-Row $ws 55  @(20000, 0, 20000, 0, 20000, $null, -20630)

Set-Row $ws 61  @(2224.5, 2119.2, 2400, 2119.2, 2400, -1907.2, -2824)

Set-Row $ws 62  @(0, 0, 0, 0, 0, $null, $null)
Clear-Cell $ws 62 14

Set-Row $ws 65  @(0, 0, 0, 0, 0, $null, $null)
Clear-Cell $ws 65 14

Set-Row $ws 74  @(246144.58, 2100.359, 5005007, 2100.359, 5005007, -1226.359, -5006755)
Set-Row $ws 77  @(246144.58, 2100.359, 5005007, 10501.795, 25025035, -6133.795, -25033771)

Set-Row $ws 80  @(28957.143, 0, 28957.143, 0, 28957.143, $null, -30953.143)
Clear-Cell $ws 80 13

Set-Row $ws 83  @(28957.143, 0, 28957.143, 0, 86871.429, $null, -96855.429)
Clear-Cell $ws 83 13

Set-Row $ws 136 @(2224.5, 2119.2, 2400, 6357.599999999999, 7200, -3807.599999999999, -12300)

# ---------------- BSM ----------------
$ws = $wb.Worksheets.Item("BSM")

Set-Row $ws 82 @(13874.728, 3595.6924, 28722.223, 3595.6924, 28722.223, -3212.6924, -29488.223)
Set-Row $ws 85 @(13874.728, 3595.6924, 28722.223, 3595.6924, 28722.223, -2269.6924, -31374.223)
Set-Row $ws 97 @(8571, 3885.2, 32000, 3885.2, 32000, -2894.2, -33982)

# ---------------- CRP ----------------
$ws = $wb.Worksheets.Item("CRP")

Set-Row $ws 31 @(2592.6333, 1374.7142, 3658.3125, 1374.7142, 3658.3125, -1079.7142, -4248.3125)
Set-Row $ws 34 @(2592.6333, 1374.7142, 3658.3125, 1374.7142, 3658.3125, -1172.7142, -4062.3125)

# ---------------- CUL ----------------
$ws = $wb.Worksheets.Item("CUL")

Set-Row $ws 127 @(1248.1666, 0, 1248.1666, 0, 3744.4998, $null, -13664.4998)

# ---------------- GSM ----------------
$ws = $wb.Worksheets.Item("GSM")

Set-Row $ws 102 @(1330.25, 1329.1666, 1340, 1329.1666, 1340, 292.8334, -4584)

# ---------------- LTW ----------------
$ws = $wb.Worksheets.Item("LTW")

Set-Row $ws 22  @(2084268.1, 3030983, 1495.6, 3030983, 1495.6, -3030688, -2085.6)
Set-Row $ws 27  @(2084268.1, 3030983, 1495.6, 3030983, 1495.6, -3030876, -1709.6)
Set-Row $ws 40  @(2550.75, 2262.2, 3031.6667, 2262.2, 3031.6667, -2126.2, -3303.6667)
Set-Row $ws 55  @(235.24445, 240.78572, 226.11765, 240.78572, 226.11765, -67.78572, -572.11765)

Set-Row $ws 62  @(30000, 0, 30000, 0, 30000, $null, -31248)
Set-Row $ws 65  @(30000, 0, 30000, 0, 90000, $null, -96240)
Set-Row $ws 127 @(39833.332, 0, 39833.332, 0, 39833.332, $null, -49753.332)
